$d = $word.ActiveDocument

# The target paragraph originally reads (spread across 3 runs, with the
# word "possible" wrapped in proofErr gramStart/gramEnd grammar-check
# markers):
#   "1.count the number of files and folder present in the directory.  if "
#   "possible"
#   " take the directory path from user."
#
# Target layout: same overall text, but re-split as:
#   "1.count the number of files and"
#   " "
#   "folder present in the directory.  if possible take the directory path from user."

$oldText = "1.count the number of files and folder present in the directory.  if possible take the directory path from user."

# Re-typing the whole sentence via Find/Replace collapses it back into a
# single run and drops the now-stale proofErr grammar markers.
$null = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $oldText, 2)

# Locate that paragraph again (now a single consolidated run).
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.StartsWith("1.count the number of files and")) {
        $target = $para
    }
}

$pStart = $target.Range.Start

$firstText = "1.count the number of files and"
$secondText = " "
$firstLen = $firstText.Length
$secondLen = $firstLen + $secondText.Length

# Re-assigning a sub-range's FormattedText to itself forces Word to split
# the run at that boundary without altering any character formatting.
$r1 = $d.Range($pStart, $pStart + $firstLen)
$r1.FormattedText = $r1.FormattedText

$r2 = $d.Range($pStart, $pStart + $secondLen)
$r2.FormattedText = $r2.FormattedText
